# H01R0 BOM: swap the JP1 single-position header for the TE Connectivity
# "5-146280-1" part (layout/description cleanup), and tidy a couple of
# leftover formatting inconsistencies elsewhere in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("H01R0")

# --- Row 9 (JP1): new description / part number / Octopart link ---------
$ws.Range("B9").Value = "Headers & Wire Housings Unshrouded 1 POS T/H"
$ws.Range("D9").Value = "5-146280-1"

$newUrl = "https://octopart.com/5-146280-1-te+connectivity+%2F+amp-40259676?r=sp"
$ws.Range("E9").Value = $newUrl
$ws.Hyperlinks.Add($ws.Range("E9"), $newUrl) | Out-Null

# Re-apply the same look as the other Octopart link cells (the Hyperlink
# insert above nudges the cell's own formatting) and wrap the part number
# the same way the rest of the multi-line rows do.
$ws.Range("E10").Copy()
$ws.Range("E9").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 27.6

# --- Row 23 / 25: align formatting + height with the rest of the table ---
$ws.Rows.Item(23).RowHeight = 27.6
$ws.Range("E10").Copy()
$ws.Range("E23").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("E10").Copy()
$ws.Range("E25").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

# --- Leftover UI state from the last manual save --------------------------
$ws.Range("A14").Select()
